$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mark execution mode ("Manual") for each test case row (2-7) in column E
$ws.Range("E2:E7").Value = "Manual"

# Update the active selection on the sheet to E8 (mirrors the saved view state)
$ws.Range("E8").Select()
